$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl2"
$ws.Range("C2").Value = "Cxcr3"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 13.34512433333333
$ws.Range("H2").Value = 40.035373
$ws.Range("I2").Value = 0.1597931014925399
$ws.Range("J2").Value = 0.1597931014925399
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.05194133333333333
$ws.Range("N2").Value = 0.155824
$ws.Range("O2").Value = 0.03788844568234288
$ws.Range("P2").Value = 0.03788844568234288
$ws.Range("Q2").Value = 0.6931635513724445
$ws.Range("R2").Value = 6.238471962352
$ws.Range("S2").Value = 0.0060543122463132
$ws.Range("T2").Value = 0.006054312246313201
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ccl2"
$ws.Range("C3").Value = "Cxcr3"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 13.34512433333333
$ws.Range("H3").Value = 40.035373
$ws.Range("I3").Value = 0.1597931014925399
$ws.Range("J3").Value = 0.1597931014925399
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.318960333333333
$ws.Range("N3").Value = 3.956881
$ws.Range("O3").Value = 0.962111554317657
$ws.Range("P3").Value = 0.9621115543176572
$ws.Range("Q3").Value = 17.60168963906811
$ws.Range("R3").Value = 158.415206751613
$ws.Range("S3").Value = 0.1537387892462267
$ws.Range("T3").Value = 0.1537387892462267
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ccl2"
$ws.Range("C4").Value = "Cxcr3"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 26.15942766666667
$ws.Range("H4").Value = 78.478283
$ws.Range("I4").Value = 0.3132302087051685
$ws.Range("J4").Value = 0.3132302087051685
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05194133333333333
$ws.Range("N4").Value = 0.155824
$ws.Range("O4").Value = 0.03788844568234288
$ws.Range("P4").Value = 0.03788844568234288
$ws.Range("Q4").Value = 1.358755552243556
$ws.Range("R4").Value = 12.228799970192
$ws.Range("S4").Value = 0.0118678057485947
$ws.Range("T4").Value = 0.0118678057485947
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ccl2"
$ws.Range("C5").Value = "Cxcr3"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 26.15942766666667
$ws.Range("H5").Value = 78.478283
$ws.Range("I5").Value = 0.3132302087051685
$ws.Range("J5").Value = 0.3132302087051685
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.318960333333333
$ws.Range("N5").Value = 3.956881
$ws.Range("O5").Value = 0.962111554317657
$ws.Range("P5").Value = 0.9621115543176572
$ws.Range("Q5").Value = 34.50324743503589
$ws.Range("R5").Value = 310.529226915323
$ws.Range("S5").Value = 0.3013624029565738
$ws.Range("T5").Value = 0.3013624029565738
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Ccl2"
$ws.Range("C6").Value = "Cxcr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 40.356388
$ws.Range("H6").Value = 121.069164
$ws.Range("I6").Value = 0.4832231039952832
$ws.Range("J6").Value = 0.4832231039952832
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.05194133333333333
$ws.Range("N6").Value = 0.155824
$ws.Range("O6").Value = 0.03788844568234288
$ws.Range("P6").Value = 0.03788844568234288
$ws.Range("Q6").Value = 2.096164601237334
$ws.Range("R6").Value = 18.865481411136
$ws.Range("S6").Value = 0.01830857232817841
$ws.Range("T6").Value = 0.01830857232817841
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Ccl2"
$ws.Range("C7").Value = "Cxcr3"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 40.356388
$ws.Range("H7").Value = 121.069164
$ws.Range("I7").Value = 0.4832231039952832
$ws.Range("J7").Value = 0.4832231039952832
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.318960333333333
$ws.Range("N7").Value = 3.956881
$ws.Range("O7").Value = 0.962111554317657
$ws.Range("P7").Value = 0.9621115543176572
$ws.Range("Q7").Value = 53.22847496860933
$ws.Range("R7").Value = 479.056274717484
$ws.Range("S7").Value = 0.4649145316671048
$ws.Range("T7").Value = 0.4649145316671048
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Ccl2"
$ws.Range("C8").Value = "Cxcr3"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.654081666666666
$ws.Range("H8").Value = 10.962245
$ws.Range("I8").Value = 0.04375358580700841
$ws.Range("J8").Value = 0.04375358580700841
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.05194133333333333
$ws.Range("N8").Value = 0.155824
$ws.Range("O8").Value = 0.03788844568234288
$ws.Range("P8").Value = 0.03788844568234288
$ws.Range("Q8").Value = 0.1897978738755556
$ws.Range("R8").Value = 1.70818086488
$ws.Range("S8").Value = 0.001657755359256567
$ws.Range("T8").Value = 0.001657755359256567
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Ccl2"
$ws.Range("C9").Value = "Cxcr3"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.654081666666666
$ws.Range("H9").Value = 10.962245
$ws.Range("I9").Value = 0.04375358580700841
$ws.Range("J9").Value = 0.04375358580700841
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.318960333333333
$ws.Range("N9").Value = 3.956881
$ws.Range("O9").Value = 0.962111554317657
$ws.Range("P9").Value = 0.9621115543176572
$ws.Range("Q9").Value = 4.819588773093889
$ws.Range("R9").Value = 43.376298957845
$ws.Range("S9").Value = 0.04209583044775184
$ws.Range("T9").Value = 0.04209583044775184